# Apply the diff: add new simulation results (columns V:AO) mirroring the
# angle values in row 2 and new binary data in row 3, then update the
# selection to reflect where the user ended up after adding the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (angle values, same repeating pattern as B2:U2) for columns V:AO
$row2Vals = @(
    0.31415926535897898,
    0.62831853071795896,
    0.94247779607693805,
    1.2566370614359199,
    1.5707963267949001,
    1.8849555921538801,
    2.1991148575128601,
    2.5132741228718301,
    2.8274333882308098,
    3.14159265358979,
    3.4557519189487702,
    3.76991118430775,
    4.0840704496667302,
    4.3982297150257104,
    4.7123889803846897,
    5.0265482457436699,
    5.3407075111026501,
    5.6548667764616303,
    5.9690260418206096,
    6.2831853071795898
)

# Row 3 (new binary simulation data) for columns V:AO
$row3Vals = @(1,0,0,1,0,1,0,0,1,1,1,0,0,1,0,0,1,0,1,1)

# Columns V (22) through AO (41)
$startCol = 22

for ($i = 0; $i -lt $row2Vals.Length; $i++) {
    $col = $startCol + $i
    $ws.Cells.Item(2, $col).Value = $row2Vals[$i]
    $ws.Cells.Item(3, $col).Value = $row3Vals[$i]
}

# Update view state: scroll/select to where the edit left off
$ws.Range("AC10").Select()
